# Weekly price update: a new record (week of 2022-12-23) is inserted as a
# new row before the current row 197, shifting all following rows down by
# one (197-213 become 198-214). The new row keeps the same fixed
# "Mercado/Categoria" columns as the rest of the sheet and carries its own
# price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 197, pushing existing rows 197:213 down to 198:214.
$ws.Rows("197:197").Insert()

# Populate the newly inserted row 197 with the new weekly record.
$ws.Range("A197").Value = 3
$ws.Range("B197").Value = "Femacal de La Calera"
$ws.Range("C197").Value = "Coquimbo"
$ws.Range("D197").Value = 44918
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = 100112030
$ws.Range("G197").Value = "Poroto granado"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 73
$ws.Range("K197").Value = 27000
$ws.Range("L197").Value = 28000
$ws.Range("M197").Value = 27479
$ws.Range("N197").Value = "$/saco 25 kilos"
$ws.Range("O197").Value = "Provincia de Limarí"
$ws.Range("P197").Value = 1099
$ws.Range("Q197").Value = 25
$ws.Range("R197").Value = "Hortaliza"
